# excel_writer: also include totals for the balance columns
#
# Previously the "Total" row on the "Gesamtergebnis" sheet left the
# balance columns (Startguthaben / Endsaldo) as "N/A" text. Now they
# should be included in the totals like the other numeric columns, so
# write 0 numeric values into C3/D3 instead of the "N/A" string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

# Replace the "N/A" placeholders in the Total row's balance columns
# (Startguthaben = C, Endsaldo = D) with numeric totals, matching the
# other summed columns on that row.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# The Total row now auto-sizes slightly differently once it holds
# numbers instead of text in every column.
$ws.Rows.Item(3).RowHeight = 13.8

# Update the selection to match the new authoring position.
$ws.Range("C6").Select()
